# Daily attendance processing - 2026-02-01 10:40:17
# Swap "System, <email>" -> "<email>, System" in the "Recorded By" column (G)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

$oldText = "System, dnasr281@gmail.com"
$newText = "dnasr281@gmail.com, System"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    if ($cell.Value2 -eq $oldText) {
        $cell.Value2 = $newText
    }
}
